$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 28
$ws.Range("I11").Value = 28
$ws.Range("K11").Value = 28
$ws.Range("M11").Value = 112

$ws.Range("H33").Value = 241.72223
$ws.Range("I33").Value = 241.72223
$ws.Range("K33").Value = 241.72223
$ws.Range("M33").Value = -12.72223

$ws.Range("H38").Value = 2351.7778
$ws.Range("I38").Value = 1531.7142
$ws.Range("J38").Value = 5222
$ws.Range("K38").Value = 4595.142599999999
$ws.Range("L38").Value = 15666
$ws.Range("M38").Value = -4223.142599999999
$ws.Range("N38").Value = -16410

$ws.Range("H39").Value = 54
$ws.Range("I39").Value = 54
$ws.Range("K39").Value = 162
$ws.Range("M39").Value = 134

$ws.Range("H53").Value = 1256.2307
$ws.Range("I53").Value = 1236.3334
$ws.Range("K53").Value = 1236.3334
$ws.Range("M53").Value = -599.3334

$ws.Range("H98").Value = 500
$ws.Range("I98").Value = 500
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 500
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 998
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 500
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 950
$ws.Range("N122").ClearContents()

$ws.Range("H127").Value = 4065.6667
$ws.Range("I127").Value = 4065.6667
$ws.Range("K127").Value = 12197.0001
$ws.Range("M127").Value = -7237.000100000001

$ws.Range("H138").Value = 2363.4285
$ws.Range("I138").Value = 1726.8572
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 5180.571599999999
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -40.57159999999931
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 10017500
$ws.Range("J24").Value = 10017500
$ws.Range("L24").Value = 10017500
$ws.Range("N24").Value = -10018248

$ws.Range("H32").Value = 6451.067
$ws.Range("I32").Value = 6451.067
$ws.Range("K32").Value = 6451.067
$ws.Range("M32").Value = -6164.067

$ws.Range("H45").Value = 3267
$ws.Range("I45").Value = 3267
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3267
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2890
$ws.Range("N45").ClearContents()

$ws.Range("H46").Value = 27408.2
$ws.Range("J46").Value = 50145.5
$ws.Range("L46").Value = 50145.5
$ws.Range("N46").Value = -50783.5

$ws.Range("H74").Value = 2471.5789
$ws.Range("I74").Value = 1428.6154
$ws.Range("K74").Value = 1428.6154
$ws.Range("M74").Value = -554.6153999999999

$ws.Range("H77").Value = 2471.5789
$ws.Range("I77").Value = 1428.6154
$ws.Range("K77").Value = 7143.076999999999
$ws.Range("M77").Value = -2775.076999999999

$ws.Range("H100").Value = 10017500
$ws.Range("J100").Value = 10017500
$ws.Range("L100").Value = 10017500
$ws.Range("N100").Value = -10019664

$ws.Range("H122").Value = 1600.25
$ws.Range("I122").Value = 1600.25
$ws.Range("K122").Value = 4800.75
$ws.Range("M122").Value = -2350.75

$ws.Range("H132").Value = 5897.8335
$ws.Range("I132").Value = 2995.6667
$ws.Range("K132").Value = 8987.000100000001
$ws.Range("M132").Value = -6457.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2931.3125
$ws.Range("I86").Value = 3245.5454
$ws.Range("J86").Value = 2240
$ws.Range("K86").Value = 3245.5454
$ws.Range("L86").Value = 2240
$ws.Range("M86").Value = -2122.5454
$ws.Range("N86").Value = -4486

$ws.Range("H89").Value = 2931.3125
$ws.Range("I89").Value = 3245.5454
$ws.Range("J89").Value = 2240
$ws.Range("K89").Value = 16227.727
$ws.Range("L89").Value = 11200
$ws.Range("M89").Value = -10611.727
$ws.Range("N89").Value = -22432

$ws.Range("H107").Value = 935.2
$ws.Range("I107").Value = 869
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 869
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 1051
$ws.Range("N107").Value = -5040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3489.2856
$ws.Range("I16").Value = 2685
$ws.Range("K16").Value = 2685
$ws.Range("M16").Value = -2398

$ws.Range("H22").Value = 536
$ws.Range("I22").Value = 385
$ws.Range("J22").Value = 2499
$ws.Range("K22").Value = 385
$ws.Range("L22").Value = 2499
$ws.Range("M22").Value = -35
$ws.Range("N22").Value = -3199

$ws.Range("H31").Value = 2115.6428
$ws.Range("I31").Value = 814.55554
$ws.Range("J31").Value = 4457.6
$ws.Range("K31").Value = 814.55554
$ws.Range("L31").Value = 4457.6
$ws.Range("M31").Value = -519.55554
$ws.Range("N31").Value = -5047.6

$ws.Range("H34").Value = 2115.6428
$ws.Range("I34").Value = 814.55554
$ws.Range("J34").Value = 4457.6
$ws.Range("K34").Value = 814.55554
$ws.Range("L34").Value = 4457.6
$ws.Range("M34").Value = -612.55554
$ws.Range("N34").Value = -4861.6

$ws.Range("H58").Value = 1791.1818
$ws.Range("I58").Value = 773.75
$ws.Range("K58").Value = 773.75
$ws.Range("M58").Value = -570.75

$ws.Range("H113").Value = 3489.2856
$ws.Range("I113").Value = 2685
$ws.Range("K113").Value = 2685
$ws.Range("M113").Value = -515

$ws.Range("H132").Value = 2650.3572
$ws.Range("I132").Value = 2757.2
$ws.Range("K132").Value = 8271.599999999999
$ws.Range("M132").Value = -5741.599999999999

$ws.Range("H134").Value = 735.1429000000001
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 1791.1818
$ws.Range("I136").Value = 773.75
$ws.Range("K136").Value = 2321.25
$ws.Range("M136").Value = 228.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2224
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2224
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2224
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -2450

$ws.Range("H19").Value = 6561.625
$ws.Range("I19").Value = 873
$ws.Range("J19").Value = 8457.833000000001
$ws.Range("K19").Value = 873
$ws.Range("L19").Value = 8457.833000000001
$ws.Range("M19").Value = -585
$ws.Range("N19").Value = -9033.833000000001

$ws.Range("H80").Value = 10966.444
$ws.Range("I80").Value = 2299.6667
$ws.Range("K80").Value = 2299.6667
$ws.Range("M80").Value = -1301.6667

$ws.Range("H83").Value = 10966.444
$ws.Range("I83").Value = 2299.6667
$ws.Range("K83").Value = 11498.3335
$ws.Range("M83").Value = -6506.333500000001

$ws.Range("H97").Value = 662.7778
$ws.Range("I97").Value = 558.125
$ws.Range("K97").Value = 558.125
$ws.Range("M97").Value = -62.125

$ws.Range("H102").Value = 1022.3333
$ws.Range("I102").Value = 1007
$ws.Range("K102").Value = 1007
$ws.Range("M102").Value = 615

$ws.Range("H113").Value = 500
$ws.Range("I113").Value = 500
$ws.Range("K113").Value = 500
$ws.Range("M113").Value = 1670

$ws.Range("H132").Value = 826
$ws.Range("I132").Value = 826
$ws.Range("K132").Value = 2478
$ws.Range("M132").Value = 52

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4733.3335
$ws.Range("I16").Value = 4825
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 4825
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -4655
$ws.Range("N16").Value = -4340

$ws.Range("H68").Value = 699
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 699
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 1846.2
$ws.Range("I82").Value = 1952.4
$ws.Range("K82").Value = 1952.4
$ws.Range("M82").Value = -1591.4

$ws.Range("H85").Value = 1846.2
$ws.Range("I85").Value = 1952.4
$ws.Range("K85").Value = 1952.4
$ws.Range("M85").Value = -704.4000000000001

$ws.Range("H123").Value = 34500
$ws.Range("J123").Value = 34500
$ws.Range("L123").Value = 34500
$ws.Range("N123").Value = -44300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1000
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1336

$ws.Range("H69").Value = 27241.5
$ws.Range("J69").Value = 27241.5
$ws.Range("L69").Value = 27241.5
$ws.Range("N69").Value = -28739.5

$ws.Range("H72").Value = 27241.5
$ws.Range("J72").Value = 27241.5
$ws.Range("L72").Value = 81724.5
$ws.Range("N72").Value = -89212.5

$ws.Range("H122").Value = 2004
$ws.Range("I122").Value = 2004
$ws.Range("K122").Value = 6012
$ws.Range("M122").Value = -3562

$ws.Range("H127").Value = 50000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 50000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 50000
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -59920

$ws.Range("H132").Value = 1743.75
$ws.Range("I132").Value = 1711.4
$ws.Range("J132").Value = 1797.6666
$ws.Range("K132").Value = 5134.200000000001
$ws.Range("L132").Value = 5392.9998
$ws.Range("M132").Value = -2604.200000000001
$ws.Range("N132").Value = -10452.9998
